{"js": "// Add the \"getAllNew\" / \"getProductNew\" / \"getLastestNew\" / \"getNumberNew\"\n// stored-procedure text block to the document, right after the\n// \"getProductByName\" procedure's closing \"END\" paragraph (i.e. immediately\n// before the final \"GO\" paragraph that ends the document).\n\n// The lines of the new block, in document order. An empty string stands\n// for a blank spacer paragraph (matches the existing \"<w:p/>\" spacer\n// paragraphs already used between procedures in this document).\nconst newLines = [\n  \"END\",\n  \"GO\",\n  \"\",\n  \"--News\",\n  \"ALTER PROCEDURE getAllNew\",\n  \"@isAdmin BIT = 0\",\n  \"AS\",\n  \"BEGIN\",\n  \"IF(@isAdmin = 0)\",\n  \"SELECT * FROM News WHERE Status = 1 ORDER BY Created DESC\",\n  \"ELSE\",\n  \"SELECT * FROM News ORDER BY Created DESC\",\n  \"END\",\n  \"GO\",\n  \"\",\n  \"ALTER PROCEDURE getProductNew\",\n  \"@proId CHAR(7),\",\n  \"@isAdmin BIT = 0\",\n  \"AS\",\n  \"BEGIN\",\n  \"IF(@isAdmin = 0)\",\n  \"SELECT * FROM News WHERE Status = 1 AND ProductId = @proId ORDER BY Created DESC\",\n  \"ELSE\",\n  \"SELECT * FROM News WHERE ProductId = @proId ORDER BY Created DESC\",\n  \"END\",\n  \"GO\",\n  \"\",\n  \"CREATE PROCEDURE getLastestNew\",\n  \"@numNew INT = NULL\",\n  \"AS\",\n  \"BEGIN\",\n  \"IF (@numNew IS NOT NULL)\",\n  \"SELECT TOP(@numNew)* FROM News WHERE Status = 1 ORDER BY Created DESC\",\n  \"ELSE\",\n  \"SELECT * FROM News WHERE Status = 1 ORDER BY Created DESC\",\n  \"END\",\n  \"GO\",\n  \"\",\n  \"CREATE PROCEDURE getNumberNew\",\n  \"AS\",\n  \"BEGIN\",\n  \"SELECT COUNT(*) FROM News\",\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that contains the \"getProductByName\" procedure's\n// SELECT statement; the procedure's \"END\" immediately follows it.\nlet selectParaIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (\n    paragraphs.items[i].text.indexOf(\n      \"SELECT * FROM Product WHERE ProductName = @name\"\n    ) !== -1\n  ) {\n    selectParaIndex = i;\n    break;\n  }\n}\n\nif (selectParaIndex === -1) {\n  throw new Error(\n    \"Could not find the getProductByName SELECT paragraph to anchor the insertion.\"\n  );\n}\n\n// Anchor on the existing \"END\" paragraph right after it; the new block is\n// inserted after that \"END\" paragraph, pushing the document's final \"GO\"\n// paragraph further down (unchanged).\nlet anchor = paragraphs.items[selectParaIndex + 1];\n\nfor (const line of newLines) {\n  anchor = anchor.insertParagraph(line, Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "# Add the \"getAllNew\" / \"getProductNew\" / \"getLastestNew\" / \"getNumberNew\"\n# stored-procedure text block to the document, right after the\n# \"getProductByName\" procedure's closing \"END\" paragraph (i.e. immediately\n# before the final \"GO\" paragraph that ends the document).\n\n$d = $word.ActiveDocument\n\n# The lines of the new block, in document order. An empty string stands\n# for a blank spacer paragraph (matches the existing spacer paragraphs\n# already used between procedures in this document).\n$newLines = @(\n  'END',\n  'GO',\n  '',\n  '--News',\n  'ALTER PROCEDURE getAllNew',\n  '@isAdmin BIT = 0',\n  'AS',\n  'BEGIN',\n  'IF(@isAdmin = 0)',\n  'SELECT * FROM News WHERE Status = 1 ORDER BY Created DESC',\n  'ELSE',\n  'SELECT * FROM News ORDER BY Created DESC',\n  'END',\n  'GO',\n  '',\n  'ALTER PROCEDURE getProductNew',\n  '@proId CHAR(7),',\n  '@isAdmin BIT = 0',\n  'AS',\n  'BEGIN',\n  'IF(@isAdmin = 0)',\n  'SELECT * FROM News WHERE Status = 1 AND ProductId = @proId ORDER BY Created DESC',\n  'ELSE',\n  'SELECT * FROM News WHERE ProductId = @proId ORDER BY Created DESC',\n  'END',\n  'GO',\n  '',\n  'CREATE PROCEDURE getLastestNew',\n  '@numNew INT = NULL',\n  'AS',\n  'BEGIN',\n  'IF (@numNew IS NOT NULL)',\n  'SELECT TOP(@numNew)* FROM News WHERE Status = 1 ORDER BY Created DESC',\n  'ELSE',\n  'SELECT * FROM News WHERE Status = 1 ORDER BY Created DESC',\n  'END',\n  'GO',\n  '',\n  'CREATE PROCEDURE getNumberNew',\n  'AS',\n  'BEGIN',\n  'SELECT COUNT(*) FROM News'\n)\n\n# Locate the paragraph that contains the \"getProductByName\" procedure's\n# SELECT statement; the procedure's \"END\" immediately follows it.\n$rng = $d.Content\n$found = $rng.Find.Execute(\"SELECT * FROM Product WHERE ProductName = @name\")\nif (-not $found) {\n    throw \"Could not find the getProductByName SELECT paragraph to anchor the insertion.\"\n}\n$selectPara = $rng.Paragraphs(1)\n$endPara = $selectPara.Next()\n\n# Insert each new line as its own paragraph right after the \"END\"\n# paragraph, pushing the document's final \"GO\" paragraph further down\n# (unchanged). Re-fetching the freshly created paragraph via .Next() after\n# each InsertParagraphAfter() keeps the anchor accurate across the loop.\n$curPara = $endPara\nforeach ($line in $newLines) {\n    $r = $curPara.Range\n    $r.Collapse(0)\n    $r.InsertParagraphAfter()\n    $curPara = $curPara.Next()\n    if ($line -ne \"\") {\n        $curPara.Range.InsertAfter($line)\n    }\n}\n"}
